$p = $ppt.ActivePresentation

# Slide 1: Title "First slide" -> split trailing space of "First " into its own run
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(6, 1).Text = " "

# Slide 3: Title "Third slide" -> split trailing space of "Third " into its own run
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(6, 1).Text = " "
